# Add data for 2023-12-06
# Updates YTD violent-crime counts across multiple sheets to reflect one additional day of data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B2").Value = 42
$ws.Range("F2").Value = 92
$ws.Range("J2").Value = 122
$ws.Range("B3").Value = 75
$ws.Range("D3").Value = 134
$ws.Range("E3").Value = 145
$ws.Range("H3").Value = 154
$ws.Range("J3").Value = 232
$ws.Range("F6").Value = 533
$ws.Range("H6").Value = 439
$ws.Range("J6").Value = 417
$ws.Range("B7").Value = 502
$ws.Range("D7").Value = 645
$ws.Range("E7").Value = 696
$ws.Range("F7").Value = 770
$ws.Range("H7").Value = 714
$ws.Range("J7").Value = 792

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J8").Value = 45
$ws.Range("H10").Value = 6
$ws.Range("B19").Value = 7
$ws.Range("B32").Value = 16
$ws.Range("J32").Value = 47
$ws.Range("E53").Value = 83
$ws.Range("F65").Value = 38
$ws.Range("D70").Value = 9
$ws.Range("H70").Value = 15
$ws.Range("J76").Value = 15
$ws.Range("F77").Value = 22
$ws.Range("B80").Value = 16
$ws.Range("F95").Value = 4
$ws.Range("J95").Value = 4
$ws.Range("B98").Value = 502
$ws.Range("D98").Value = 645
$ws.Range("E98").Value = 696
$ws.Range("F98").Value = 770
$ws.Range("H98").Value = 714
$ws.Range("J98").Value = 792

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F6").Value = 11
$ws.Range("F7").Value = 22

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J5").Value = 25
$ws.Range("J6").Value = 45

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("B2").Value = 1
$ws.Range("J2").Value = 3
$ws.Range("B7").Value = 16
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("B2").Value = 1
$ws.Range("B5").Value = 16

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E3").Value = 14
$ws.Range("E7").Value = 83

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("F5").Value = 31
$ws.Range("F6").Value = 38

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("H5").Value = 5
$ws.Range("H6").Value = 6

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("B3").Value = 2
$ws.Range("B6").Value = 7

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("D3").Value = 3
$ws.Range("H3").Value = 4
$ws.Range("D5").Value = 9
$ws.Range("H5").Value = 15

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("E2").Value = 1
$ws.Range("H3").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("H5").Value = 4
